$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# stay as text (matching the source data which stores them as strings),
# otherwise Excel auto-converts them to numeric values on assignment.
$textCells = @("D8","D10","D13","D15","D16","D19","D20","D22","D25","D26","D29","D31","D33","D34","D42","D46","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (prices and 1h volume deltas).
$ws.Range('D2').Value = '37.284.18'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.062.74'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '56.58'
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').Value = '0.0759'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '2.366.74'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '14.58'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').Value = '0.774'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '5.11'
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('D17').Value = '2.063.36'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '37.241.48'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = '6.31'
$ws.Range('E19').Value = '  +6.06%  '
$ws.Range('D20').Value = '69.44'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').Value = '226.10'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('D26').Value = '166.18'
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('E27').Value = '  +4.25%  '
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('D29').Value = '19.03'
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('D31').Value = '0.118'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').Value = '0.0614'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +3.81%  '
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('E38').Value = '  -3.27%  '
$ws.Range('E39').Value = '  -4.68%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').Value = '1.463.63'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').Value = '96.06'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('E43').Value = '  -3.07%  '
$ws.Range('E44').Value = '  +1.71%  '
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('D46').Value = '4.21'
$ws.Range('E46').Value = '  -5.25%  '
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '7.14'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '14.97'
$ws.Range('E49').Value = '  -7.08%  '
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '2.252.99'
$ws.Range('E51').Value = '  -0.26%  '
